$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2025-08-19 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-20 Wednesday", 2) | Out-Null

# Update each cell of the answer table, in row-major order, so that
# duplicate "old" values (e.g. "67-48=19") are replaced independently
# rather than via a single document-wide Find/Replace.
$t = $d.Tables(1)
$cols = 5
$pairs = @(
    @("4+37=41", "27+27=54"),
    @("61-25=36", "25-8=17"),
    @("54-15=39", "7+7=14"),
    @("29+46=75", "77-58=19"),
    @("63+8=71", "11-6=5"),
    @("72-17=55", "39+25=64"),
    @("83-59=24", "4+88=92"),
    @("71-36=35", "34-26=8"),
    @("60-57=3", "56-48=8"),
    @("24+28=52", "70-8=62"),
    @("70-18=52", "54+8=62"),
    @("14+29=43", "17+65=82"),
    @("28+49=77", "72-64=8"),
    @("82-75=7", "77+16=93"),
    @("8+84=92", "90-23=67"),
    @("83-68=15", "49+26=75"),
    @("93-19=74", "93-29=64"),
    @("39+8=47", "92-67=25"),
    @("40-33=7", "12+69=81"),
    @("45+26=71", "58+8=66"),
    @("62-59=3", "73-28=45"),
    @("38+47=85", "81-68=13"),
    @("32-19=13", "9+12=21"),
    @("30-9=21", "86-67=19"),
    @("74-68=6", "59+13=72"),
    @("73-44=29", "38+57=95"),
    @("28+27=55", "95-79=16"),
    @("3+79=82", "57-49=8"),
    @("32+19=51", "46-7=39"),
    @("94-56=38", "88-39=49"),
    @("6+5=11", "84+8=92"),
    @("83-27=56", "20-5=15"),
    @("48+17=65", "70-43=27"),
    @("29+69=98", "32-6=26"),
    @("68-9=59", "55-8=47"),
    @("90-46=44", "16+77=93"),
    @("70-34=36", "93-47=46"),
    @("47+29=76", "97-88=9"),
    @("7+8=15", "60-32=28"),
    @("26+66=92", "8+54=62"),
    @("2+79=81", "65+18=83"),
    @("82-59=23", "19+67=86"),
    @("75-28=47", "21-18=3"),
    @("98-79=19", "88-69=19"),
    @("15+76=91", "62+9=71"),
    @("60-53=7", "29+33=62"),
    @("80-69=11", "60-45=15"),
    @("43-7=36", "8+23=31"),
    @("18+54=72", "93-48=45"),
    @("18+18=36", "91-13=78"),
    @("91-64=27", "91-44=47"),
    @("58+15=73", "57-8=49"),
    @("27+47=74", "35+58=93"),
    @("72-18=54", "6+75=81"),
    @("38+36=74", "92-16=76"),
    @("55+37=92", "64+8=72"),
    @("48+43=91", "16+79=95"),
    @("38+25=63", "38+54=92"),
    @("30-19=11", "66-39=27"),
    @("22+59=81", "60-36=24"),
    @("70-15=55", "18+53=71"),
    @("90-22=68", "94-77=17"),
    @("54-45=9", "69+14=83"),
    @("67-48=19", "41-16=25"),
    @("39+4=43", "24-9=15"),
    @("13+59=72", "36-27=9"),
    @("70-12=58", "15+57=72"),
    @("18+63=81", "73-46=27"),
    @("67-48=19", "31-22=9"),
    @("2+59=61", "13+68=81"),
    @("43+48=91", "81-45=36"),
    @("5+49=54", "83-57=26"),
    @("68+24=92", "70-67=3"),
    @("48+36=84", "52-49=3"),
    @("86+8=94", "7+36=43"),
    @("46-29=17", "23-8=15"),
    @("82-79=3", "85-8=77"),
    @("40-2=38", "30-13=17"),
    @("80-34=46", "92-64=28"),
    @("70-9=61", "58-29=29"),
    @("15-7=8", "62-34=28"),
    @("64-38=26", "17+25=42"),
    @("15-8=7", "57+25=82"),
    @("74+7=81", "71-65=6"),
    @("48+29=77", "60-56=4"),
    @("29+17=46", "17+74=91"),
    @("45-28=17", "66-49=17"),
    @("82-29=53", "94-69=25"),
    @("23+58=81", "7+15=22"),
    @("37+25=62", "72-35=37"),
    @("26+26=52", "61-35=26"),
    @("58+15=73", "77-48=29"),
    @("38+25=63", "70-42=28"),
    @("3+59=62", "40-17=23"),
    @("52-25=27", "8+29=37"),
    @("46+17=63", "69+24=93"),
    @("41-35=6", "73+19=92"),
    @("36+15=51", "89+4=93"),
    @("58-19=39", "54+29=83"),
    @("15+58=73", "91-8=83")
)

for ($i = 0; $i -lt $pairs.Count; $i++) {
    $row = [int][Math]::Floor($i / $cols) + 1
    $col = ($i % $cols) + 1
    $old = $pairs[$i][0]
    $new = $pairs[$i][1]
    $cell = $t.Cell($row, $col)
    $cell.Range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 0, $false, $new, 1) | Out-Null
}

Write-Output "done"
